$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 3099.503889238888
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 3748.219677045132

$ws.Range("B3").Value = 0.6753301551942219
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 2.290389397800092
